$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-style the existing Chapter 12 "Working with tables" rows (62-68).
#    They previously carried bespoke xfs (numFmt "00" w/ applyFill/applyBorder,
#    and General w/ applyFill/applyBorder); normalise them back to the plain
#    styles used everywhere else in the sheet (numFmt "00" == style of C2,
#    General == default "Normal").
# ---------------------------------------------------------------------------
for ($r = 62; $r -le 68; $r++) {
    foreach ($col in @("C", "E", "G")) {
        $cell = $ws.Range("$col$r")
        $cell.Style = "Normal"
        $cell.NumberFormat = "00"
    }
    foreach ($col in @("D", "F", "H")) {
        $ws.Range("$col$r").Style = "Normal"
    }
    $iCell = $ws.Range("I$r")
    $iCell.Style = "Normal"
    $iCell.NumberFormat = "00"
}

# ---------------------------------------------------------------------------
# 2. Finish the "Using tables as filters" section: append rows 69-73.
# ---------------------------------------------------------------------------

# Row 69 - Implementing OR conditions
$ws.Range("C69").NumberFormat = "00"
$ws.Range("C69").Value2 = 12
$ws.Range("D69").Value = "Working with tables"
$ws.Range("E69").NumberFormat = "00"
$ws.Range("E69").Value2 = 3
$ws.Range("I69").NumberFormat = "00"
$ws.Range("I69").Value = "SUMMARIZE, CROSSJOIN, CALCULATE"
$ws.Range("F69").Value = "Using tables as filters"
$ws.Range("G69").NumberFormat = "00"
$ws.Range("G69").Value2 = 1
$ws.Range("H69").Value = "Implementing OR conditions"
$ws.Range("B69").Formula = "=_xlfn.CONCAT(TEXT(C69,""00""),TEXT(E69,""00""),TEXT(G69,""00""))"

# Row 70 - Narrowing sales computation to the first year's customers
$ws.Range("C70").NumberFormat = "00"
$ws.Range("C70").Value2 = 12
$ws.Range("D70").Value = "Working with tables"
$ws.Range("E70").NumberFormat = "00"
$ws.Range("E70").Value2 = 3
$ws.Range("F70").Value = "Using tables as filters"
$ws.Range("G70").NumberFormat = "00"
$ws.Range("G70").Value2 = 2
$ws.Range("H70").Value = "Narrowing sales computation to the first years customers"
$ws.Range("I70").NumberFormat = "00"
$ws.Range("I70").Value = "CALCULATETABLE, ALLSELECTED, FIRSTNONBLANK, KEEPFILTERS"
$ws.Range("B70").Formula = "=_xlfn.CONCAT(TEXT(C70,""00""),TEXT(E70,""00""),TEXT(G70,""00""))"

# Row 71 - Computing new customers
$ws.Range("C71").NumberFormat = "00"
$ws.Range("C71").Value2 = 12
$ws.Range("D71").Value = "Working with tables"
$ws.Range("E71").NumberFormat = "00"
$ws.Range("E71").Value2 = 3
$ws.Range("F71").Value = "Using tables as filters"
$ws.Range("G71").NumberFormat = "00"
$ws.Range("G71").Value2 = 3
$ws.Range("H71").Value = "Computing new customers"
$ws.Range("I71").NumberFormat = "00"
$ws.Range("I71").Value = "CALCULATETABLE, ADDCOLUMNS, VALUES, FILTER, COUNTROWS"
$ws.Range("B71").Formula = "=_xlfn.CONCAT(TEXT(C71,""00""),TEXT(E71,""00""),TEXT(G71,""00""))"

# Row 72 - Reusing table expressions with DETAILROWS
$ws.Range("C72").NumberFormat = "00"
$ws.Range("C72").Value2 = 12
$ws.Range("D72").Value = "Working with tables"
$ws.Range("E72").NumberFormat = "00"
$ws.Range("E72").Value2 = 3
$ws.Range("F72").Value = "Using tables as filters"
$ws.Range("G72").NumberFormat = "00"
$ws.Range("G72").Value2 = 4
$ws.Range("H72").Value = "Reusing table expressions with DETAILROWS"
$ws.Range("B72").Formula = "=_xlfn.CONCAT(TEXT(C72,""00""),TEXT(E72,""00""),TEXT(G72,""00""))"

# Row 73 - chapter total marker row (chapter/section only)
$ws.Range("C73").NumberFormat = "00"
$ws.Range("C73").Value2 = 12
$ws.Range("D73").Value = "Working with tables"
$ws.Range("B73").Formula = "=_xlfn.CONCAT(TEXT(C73,""00""),TEXT(E73,""00""),TEXT(G73,""00""))"

# ---------------------------------------------------------------------------
# 3. Restore the selection that was active when the edit was saved.
# ---------------------------------------------------------------------------
$ws.Range("I72").Select()
